# Applies the "all team member billable hours" update described by the commit.

$wb = $excel.ActiveWorkbook

# --- Workbook-level metadata (window size / active selection on Status Report sheet) ---
$wsStatus = $wb.Worksheets.Item("Status Report")
$wsStatus.Activate()

# New cell content / edits on the "Status Report" sheet
# Row 22 (Yiqi Wang): Role/Task reset to N/A, hours 1 -> 0
$wsStatus.Range("C22").Value = "N/A"
$wsStatus.Range("E22").Value = 0

# Row 23: new team member row - Ian Jackson, org chart + status report, 0.5 hours
$wsStatus.Range("C23").Value = "org chart + status report"
$wsStatus.Range("A23").Value = "Ian Jackson"
$wsStatus.Range("E23").Value = 0.5

# J12: add "Requirements - Identify on 10/4/23"
$wsStatus.Range("J12").Value = "Requirements - Identify on 10/4/23"

# Select the range shown active in the saved file
$wsStatus.Range("A12:E16").Select()

$wb.Save()
